# Update recomputed profit-analysis figures across the Leve profit sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR). These H:N columns
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ,
# LeveProfitNQ/HQ) are refreshed market-price-derived values produced by a
# scheduled data-refresh run; each row below is pinned to its new value.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 389
$ws.Cells.Item(33, 9).Value = 381.2353
$ws.Cells.Item(33, 10).Value = 433
$ws.Cells.Item(33, 11).Value = 381.2353
$ws.Cells.Item(33, 12).Value = 433
$ws.Cells.Item(33, 13).Value = -152.2353
$ws.Cells.Item(33, 14).Value = -891
$ws.Cells.Item(74, 8).Value = 8215.579
$ws.Cells.Item(74, 9).Value = 5624.25
$ws.Cells.Item(74, 10).Value = 8906.6
$ws.Cells.Item(74, 11).Value = 5624.25
$ws.Cells.Item(74, 12).Value = 8906.6
$ws.Cells.Item(74, 13).Value = -4688.25
$ws.Cells.Item(74, 14).Value = -10778.6
$ws.Cells.Item(77, 8).Value = 8215.579
$ws.Cells.Item(77, 9).Value = 5624.25
$ws.Cells.Item(77, 10).Value = 8906.6
$ws.Cells.Item(77, 11).Value = 28121.25
$ws.Cells.Item(77, 12).Value = 44533
$ws.Cells.Item(77, 13).Value = -23441.25
$ws.Cells.Item(77, 14).Value = -53893
$ws.Cells.Item(98, 8).Value = 1511.5883
$ws.Cells.Item(98, 9).Value = 1416.5476
$ws.Cells.Item(98, 10).Value = 1955.1111
$ws.Cells.Item(98, 11).Value = 1416.5476
$ws.Cells.Item(98, 12).Value = 1955.1111
$ws.Cells.Item(98, 13).Value = 81.4523999999999
$ws.Cells.Item(98, 14).Value = -4951.1111
$ws.Cells.Item(113, 8).Value = 2874.1667
$ws.Cells.Item(113, 9).Value = 2243.3333
$ws.Cells.Item(113, 10).Value = 4766.6665
$ws.Cells.Item(113, 11).Value = 2243.3333
$ws.Cells.Item(113, 12).Value = 4766.6665
$ws.Cells.Item(113, 13).Value = 1010.6667
$ws.Cells.Item(113, 14).Value = -11274.6665
$ws.Cells.Item(116, 8).Value = 13217
$ws.Cells.Item(116, 9).Value = 11851.368
$ws.Cells.Item(116, 10).Value = 17541.5
$ws.Cells.Item(116, 11).Value = 11851.368
$ws.Cells.Item(116, 12).Value = 17541.5
$ws.Cells.Item(116, 13).Value = -8409.368
$ws.Cells.Item(116, 14).Value = -24425.5
$ws.Cells.Item(122, 8).Value = 1511.5883
$ws.Cells.Item(122, 9).Value = 1416.5476
$ws.Cells.Item(122, 10).Value = 1955.1111
$ws.Cells.Item(122, 11).Value = 4249.642800000001
$ws.Cells.Item(122, 12).Value = 5865.3333
$ws.Cells.Item(122, 13).Value = -1799.642800000001
$ws.Cells.Item(122, 14).Value = -10765.3333
$ws.Cells.Item(125, 8).Value = 3079.9443
$ws.Cells.Item(125, 10).Value = 2754.9092
$ws.Cells.Item(125, 12).Value = 24794.1828
$ws.Cells.Item(125, 14).Value = -29714.1828
$ws.Cells.Item(132, 8).Value = 47275.727
$ws.Cells.Item(132, 9).Value = 52197
$ws.Cells.Item(132, 11).Value = 156591
$ws.Cells.Item(132, 13).Value = -154061
$ws.Cells.Item(137, 8).Value = 1055149.8
$ws.Cells.Item(137, 9).Value = 604797.0600000001
$ws.Cells.Item(137, 10).Value = 1505502.5
$ws.Cells.Item(137, 11).Value = 1814391.18
$ws.Cells.Item(137, 12).Value = 4516507.5
$ws.Cells.Item(137, 13).Value = -1811841.18
$ws.Cells.Item(137, 14).Value = -4521607.5
$ws.Cells.Item(138, 8).Value = 3722.4412
$ws.Cells.Item(138, 9).Value = 2329.1667
$ws.Cells.Item(138, 10).Value = 5289.875
$ws.Cells.Item(138, 11).Value = 6987.500100000001
$ws.Cells.Item(138, 12).Value = 15869.625
$ws.Cells.Item(138, 13).Value = -1847.500100000001
$ws.Cells.Item(138, 14).Value = -26149.625
$ws.Cells.Item(141, 8).Value = 1637.5
$ws.Cells.Item(141, 9).Value = 1275
$ws.Cells.Item(141, 11).Value = 3825
$ws.Cells.Item(141, 13).Value = 1355
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 9809273
$ws.Cells.Item(32, 9).Value = 10642402
$ws.Cells.Item(32, 11).Value = 10642402
$ws.Cells.Item(32, 13).Value = -10642115
$ws.Cells.Item(45, 8).Value = 4053.3704
$ws.Cells.Item(45, 9).Value = 3820.8823
$ws.Cells.Item(45, 11).Value = 3820.8823
$ws.Cells.Item(45, 13).Value = -3443.8823
$ws.Cells.Item(61, 8).Value = 1118813.9
$ws.Cells.Item(61, 9).Value = 1198514.9
$ws.Cells.Item(61, 11).Value = 1198514.9
$ws.Cells.Item(61, 13).Value = -1198302.9
$ws.Cells.Item(110, 8).Value = 588.8
$ws.Cells.Item(110, 9).Value = 588.8
$ws.Cells.Item(110, 11).Value = 588.8
$ws.Cells.Item(110, 13).Value = 1456.2
$ws.Cells.Item(122, 8).Value = 2576.7273
$ws.Cells.Item(122, 9).Value = 1546.0667
$ws.Cells.Item(122, 10).Value = 4785.2856
$ws.Cells.Item(122, 11).Value = 4638.2001
$ws.Cells.Item(122, 12).Value = 14355.8568
$ws.Cells.Item(122, 13).Value = -2188.2001
$ws.Cells.Item(122, 14).Value = -19255.8568
$ws.Cells.Item(132, 8).Value = 834517.1
$ws.Cells.Item(132, 9).Value = 1165019.1
$ws.Cells.Item(132, 10).Value = 8262
$ws.Cells.Item(132, 11).Value = 3495057.3
$ws.Cells.Item(132, 12).Value = 24786
$ws.Cells.Item(132, 13).Value = -3492527.3
$ws.Cells.Item(132, 14).Value = -29846
$ws.Cells.Item(136, 8).Value = 1118813.9
$ws.Cells.Item(136, 9).Value = 1198514.9
$ws.Cells.Item(136, 11).Value = 3595544.7
$ws.Cells.Item(136, 13).Value = -3592994.7
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1408.5862
$ws.Cells.Item(86, 9).Value = 1362.174
$ws.Cells.Item(86, 11).Value = 1362.174
$ws.Cells.Item(86, 13).Value = -239.174
$ws.Cells.Item(89, 8).Value = 1408.5862
$ws.Cells.Item(89, 9).Value = 1362.174
$ws.Cells.Item(89, 11).Value = 6810.87
$ws.Cells.Item(89, 13).Value = -1194.87
$ws.Cells.Item(94, 8).Value = 1951.1666
$ws.Cells.Item(94, 9).Value = 1404.3846
$ws.Cells.Item(94, 10).Value = 3372.8
$ws.Cells.Item(94, 11).Value = 1404.3846
$ws.Cells.Item(94, 12).Value = 3372.8
$ws.Cells.Item(94, 13).Value = -953.3846000000001
$ws.Cells.Item(94, 14).Value = -4274.8
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 866.9091
$ws.Cells.Item(16, 9).Value = 608.2
$ws.Cells.Item(16, 10).Value = 1082.5
$ws.Cells.Item(16, 11).Value = 608.2
$ws.Cells.Item(16, 12).Value = 1082.5
$ws.Cells.Item(16, 13).Value = -321.2
$ws.Cells.Item(16, 14).Value = -1656.5
$ws.Cells.Item(113, 8).Value = 866.9091
$ws.Cells.Item(113, 9).Value = 608.2
$ws.Cells.Item(113, 10).Value = 1082.5
$ws.Cells.Item(113, 11).Value = 608.2
$ws.Cells.Item(113, 12).Value = 1082.5
$ws.Cells.Item(113, 13).Value = 1561.8
$ws.Cells.Item(113, 14).Value = -5422.5
$ws.Cells.Item(132, 8).Value = 5823467.5
$ws.Cells.Item(132, 10).Value = 83334410
$ws.Cells.Item(132, 12).Value = 250003230
$ws.Cells.Item(132, 14).Value = -250008290
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(7, 8).Value = 3416903.2
$ws.Cells.Item(7, 9).Value = 2222510
$ws.Cells.Item(7, 11).Value = 6667530
$ws.Cells.Item(7, 13).Value = -6667418
$ws.Cells.Item(80, 8).Value = 3912.1428
$ws.Cells.Item(80, 9).Value = 3848
$ws.Cells.Item(80, 10).Value = 4072.5
$ws.Cells.Item(80, 11).Value = 11544
$ws.Cells.Item(80, 12).Value = 12217.5
$ws.Cells.Item(80, 13).Value = -10608
$ws.Cells.Item(80, 14).Value = -14089.5
$ws.Cells.Item(83, 8).Value = 3912.1428
$ws.Cells.Item(83, 9).Value = 3848
$ws.Cells.Item(83, 10).Value = 4072.5
$ws.Cells.Item(83, 11).Value = 34632
$ws.Cells.Item(83, 12).Value = 36652.5
$ws.Cells.Item(83, 13).Value = -29952
$ws.Cells.Item(83, 14).Value = -46012.5
$ws.Cells.Item(86, 8).Value = 753.125
$ws.Cells.Item(86, 9).Value = 683.3333
$ws.Cells.Item(86, 10).Value = 962.5
$ws.Cells.Item(86, 11).Value = 2049.9999
$ws.Cells.Item(86, 12).Value = 2887.5
$ws.Cells.Item(86, 13).Value = -863.9998999999998
$ws.Cells.Item(86, 14).Value = -5259.5
$ws.Cells.Item(89, 8).Value = 753.125
$ws.Cells.Item(89, 9).Value = 683.3333
$ws.Cells.Item(89, 10).Value = 962.5
$ws.Cells.Item(89, 11).Value = 6149.9997
$ws.Cells.Item(89, 12).Value = 8662.5
$ws.Cells.Item(89, 13).Value = -221.9997000000003
$ws.Cells.Item(89, 14).Value = -20518.5
$ws.Cells.Item(92, 8).Value = 1139.2858
$ws.Cells.Item(92, 9).Value = 934.9286
$ws.Cells.Item(92, 11).Value = 2804.7858
$ws.Cells.Item(92, 13).Value = -1556.7858
$ws.Cells.Item(131, 8).Value = 14212.944
$ws.Cells.Item(131, 9).Value = 999
$ws.Cells.Item(131, 10).Value = 14990.235
$ws.Cells.Item(131, 11).Value = 2997
$ws.Cells.Item(131, 12).Value = 44970.705
$ws.Cells.Item(131, 13).Value = 2043
$ws.Cells.Item(131, 14).Value = -55050.705
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2585.1035
$ws.Cells.Item(102, 9).Value = 1502.7727
$ws.Cells.Item(102, 11).Value = 1502.7727
$ws.Cells.Item(102, 13).Value = 119.2273
$ws.Cells.Item(113, 8).Value = 2284.5454
$ws.Cells.Item(113, 9).Value = 866.2308
$ws.Cells.Item(113, 10).Value = 4333.222
$ws.Cells.Item(113, 11).Value = 866.2308
$ws.Cells.Item(113, 12).Value = 4333.222
$ws.Cells.Item(113, 13).Value = 1303.7692
$ws.Cells.Item(113, 14).Value = -8673.222
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 52521.2
$ws.Cells.Item(22, 9).Value = 169278.5
$ws.Cells.Item(22, 10).Value = 2482.3572
$ws.Cells.Item(22, 11).Value = 169278.5
$ws.Cells.Item(22, 12).Value = 2482.3572
$ws.Cells.Item(22, 13).Value = -168983.5
$ws.Cells.Item(22, 14).Value = -3072.3572
$ws.Cells.Item(27, 8).Value = 52521.2
$ws.Cells.Item(27, 9).Value = 169278.5
$ws.Cells.Item(27, 10).Value = 2482.3572
$ws.Cells.Item(27, 11).Value = 169278.5
$ws.Cells.Item(27, 12).Value = 2482.3572
$ws.Cells.Item(27, 13).Value = -169171.5
$ws.Cells.Item(27, 14).Value = -2696.3572
$ws.Cells.Item(40, 8).Value = 4730.773
$ws.Cells.Item(40, 9).Value = 4059.889
$ws.Cells.Item(40, 10).Value = 7749.75
$ws.Cells.Item(40, 11).Value = 4059.889
$ws.Cells.Item(40, 12).Value = 7749.75
$ws.Cells.Item(40, 13).Value = -3923.889
$ws.Cells.Item(40, 14).Value = -8021.75
$ws.Cells.Item(46, 8).Value = 1158.3077
$ws.Cells.Item(46, 10).Value = 1165.0769
$ws.Cells.Item(46, 12).Value = 1165.0769
$ws.Cells.Item(46, 14).Value = -1541.0769
$ws.Cells.Item(68, 8).Value = 3291.8125
$ws.Cells.Item(68, 9).Value = 2879.4546
$ws.Cells.Item(68, 10).Value = 4199
$ws.Cells.Item(68, 11).Value = 2879.4546
$ws.Cells.Item(68, 12).Value = 4199
$ws.Cells.Item(68, 13).Value = -2130.4546
$ws.Cells.Item(68, 14).Value = -5697
$ws.Cells.Item(71, 8).Value = 3291.8125
$ws.Cells.Item(71, 9).Value = 2879.4546
$ws.Cells.Item(71, 10).Value = 4199
$ws.Cells.Item(71, 11).Value = 14397.273
$ws.Cells.Item(71, 12).Value = 20995
$ws.Cells.Item(71, 13).Value = -10653.273
$ws.Cells.Item(71, 14).Value = -28483
$ws.Cells.Item(136, 8).Value = 4916.25
$ws.Cells.Item(136, 10).Value = 10658.167
$ws.Cells.Item(136, 12).Value = 31974.501
$ws.Cells.Item(136, 14).Value = -37074.501
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 3208.2354
$ws.Cells.Item(107, 9).Value = 1440.8334
$ws.Cells.Item(107, 10).Value = 4172.273
$ws.Cells.Item(107, 11).Value = 4322.5002
$ws.Cells.Item(107, 12).Value = 12516.819
$ws.Cells.Item(107, 13).Value = -2402.5002
$ws.Cells.Item(107, 14).Value = -16356.819
